$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Card Library")
$ws2 = $wb.Worksheets.Item("Deck Saves")

# ---------------------------------------------------------------------------
# Seed the brand-new card-effect strings first, in authoring order, so the
# shared-string table comes out in the same sequence it was originally typed.
# ---------------------------------------------------------------------------

$ws1.Range("B8").Value = "Spider man gives you 20 speed"
$ws1.Range("B9").Value = "Gain 20 hp (Consume)"
$ws1.Range("B7").Value = "dance dance dance! Ally gains 20 speed"
$ws1.Range("B6").Value = "Heal 5 for the next 3 turns"
$ws1.Range("B4").Value = "Powerful.  Deadly. Terrifying. If hp is less than 40 set hp to 0"
$ws1.Range("B2").Value = "summon a clone of yourself with 1 hp"
$ws1.Range("B5").Value = "Attack for 25 dmg"
$ws1.Range("B3").Value = "Attack for 14 dmg"

# ---------------------------------------------------------------------------
# "Card Library" sheet: replace flavor-text descriptions with mechanical
# effect text, and re-point the artwork column where applicable.
# (Title / AP columns stay the same; only Description (B) and Art (D) move.)
# ---------------------------------------------------------------------------

$ws1.Range("A1").Value = "Title"
$ws1.Range("B1").Value = "Description"
$ws1.Range("C1").Value = "AP"
$ws1.Range("D1").Value = "Art"

$ws1.Range("A2").Value = "Ghost Clone"
$ws1.Range("B2").Value = "summon a clone of yourself with 1 hp"
$ws1.Range("C2").Value = 2
$ws1.Range("D2").Value = "ghost-ally"

$ws1.Range("A3").Value = "Stiletto"
$ws1.Range("B3").Value = "Attack for 14 dmg"
$ws1.Range("C3").Value = 1
$ws1.Range("D3").Value = "stiletto"

$ws1.Range("A4").Value = "Staff of Death"
$ws1.Range("B4").Value = "Powerful.  Deadly. Terrifying. If hp is less than 40 set hp to 0"
$ws1.Range("C4").Value = 3
$ws1.Range("D4").Value = "skull-staff"

$ws1.Range("A5").Value = "Piercing Stab"
$ws1.Range("B5").Value = "Attack for 25 dmg"
$ws1.Range("C5").Value = 1
$ws1.Range("D5").Value = "stiletto"

$ws1.Range("A6").Value = "Heal"
$ws1.Range("B6").Value = "Heal 5 for the next 3 turns"
$ws1.Range("C6").Value = 4
$ws1.Range("D6").Value = "half-heart"

$ws1.Range("A7").Value = "Dance Party"
$ws1.Range("B7").Value = "dance dance dance! Ally gains 20 speed"
$ws1.Range("C7").Value = 5
$ws1.Range("D7").Value = "fire-silhouette"

$ws1.Range("A8").Value = "SpiderMan slays"
$ws1.Range("B8").Value = "Spider man gives you 20 speed"
$ws1.Range("C8").Value = 3
$ws1.Range("D8").Value = "spider-mask"

$ws1.Range("A9").Value = "Mango"
$ws1.Range("B9").Value = "Gain 20 hp (Consume)"
$ws1.Range("C9").Value = 1
$ws1.Range("D9").Value = "peach"

$ws1.Range("A10").Value = "Flame Bash"
$ws1.Range("B10").Value = "Charge forwards in a ball of fire (Deal 10 dmg, Gain 5 Speed) "
$ws1.Range("C10").Value = 2
$ws1.Range("D10").Value = "fire-dash"

$ws1.Range("A11").Value = "Dragon Breath"
$ws1.Range("B11").Value = "Unleash a breath of fire on a target dealing 12 dmg"
$ws1.Range("C11").Value = 1
$ws1.Range("D11").Value = "fire-breath"

$ws1.Range("A12").Value = "Strawberry"
$ws1.Range("B12").Value = "Heals 10 hp (Consume)"
$ws1.Range("C12").Value = 1
$ws1.Range("D12").Value = "strawberry"

$ws1.Range("A13").Value = "Motion Sickness"
$ws1.Range("B13").Value = "Reduce target speed by 15"
$ws1.Range("C13").Value = 1
$ws1.Range("D13").Value = "vomiting"

$ws1.Range("A14").Value = "Map it Out"
$ws1.Range("B14").Value = "Gain 7 speed"
$ws1.Range("C14").Value = 1
$ws1.Range("D14").Value = "slalom"

# ---------------------------------------------------------------------------
# "Deck Saves" sheet: card-name column tracks the (now-renumbered) shared
# strings, and a few deck counts were tweaked.
# ---------------------------------------------------------------------------

$ws2.Range("A2").Value = "Ghost Clone"
$ws2.Range("B2").Value = 4
$ws2.Range("C2").Value = 1

$ws2.Range("A3").Value = "Stiletto"
$ws2.Range("B3").Value = 3
$ws2.Range("C3").Value = 3

$ws2.Range("A4").Value = "Staff of Death"
$ws2.Range("B4").Value = 2
$ws2.Range("C4").Value = 0

$ws2.Range("A5").Value = "Piercing Stab"
$ws2.Range("B5").Value = 2
$ws2.Range("C5").Value = 2

$ws2.Range("A6").Value = "Heal"
$ws2.Range("B6").Value = 1
$ws2.Range("C6").Value = 10

$ws2.Range("A7").Value = "SpiderMan slays"
$ws2.Range("B7").Value = 3
$ws2.Range("C7").Value = 1

$ws2.Range("A8").Value = "Dance Party"
$ws2.Range("B8").Value = 2
$ws2.Range("C8").Value = 5

$ws2.Range("A9").Value = "Mango"
$ws2.Range("B9").Value = 2
$ws2.Range("C9").Value = 0

$ws2.Range("A10").Value = "Flame Bash"
$ws2.Range("B10").Value = 5
$ws2.Range("C10").Value = 5

$ws2.Range("A11").Value = "Dragon Breath"
$ws2.Range("B11").Value = 2
$ws2.Range("C11").Value = 2

$ws2.Range("A12").Value = "Strawberry"
$ws2.Range("B12").Value = 2
$ws2.Range("C12").Value = 2

$ws2.Range("A13").Value = "Motion Sickness"
$ws2.Range("B13").Value = 4
$ws2.Range("C13").Value = 4

$ws2.Range("A14").Value = "Map it Out"
$ws2.Range("B14").Value = 4
$ws2.Range("C14").Value = 4

# ---------------------------------------------------------------------------
# View/selection state: active tab moves from "Card Library" to
# "Deck Saves", and each sheet's remembered selection changes.
# ---------------------------------------------------------------------------

$ws1.Range("B2").Select() | Out-Null
$ws2.Range("C22").Select() | Out-Null
